# Updated cryptos list - applies latest price/volume scrape plus a row swap (Polkadot <-> WrappedEther)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, $Text)
    # Force the cell to remain a text value (matches the source data, which is
    # stored as text even when it looks numeric, e.g. "219.97" or "0.5280").
    $r = $ws.Range($Cell)
    $r.NumberFormat = "@"
    $r.Value = $Text
    $r.Style = "Normal"
}

# Row 2
Set-TextCell 'D2' '26.507.05'
$ws.Range('E2').Value = '  +1.85%  '

# Row 3
Set-TextCell 'D3' '1.671.75'
$ws.Range('E3').Value = '  +1.56%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
Set-TextCell 'D5' '219.97'
$ws.Range('E5').Value = '  +2.26%  '

# Row 6
Set-TextCell 'D6' '0.5280'
$ws.Range('E6').Value = '  +1.23%  '

# Row 7
Set-TextCell 'D7' '1.001'
$ws.Range('E7').Value = '  -0.02%  '

# Row 8
Set-TextCell 'D8' '0.2678'
$ws.Range('E8').Value = '  +2.68%  '

# Row 9
Set-TextCell 'D9' '0.06375'
$ws.Range('E9').Value = '  +0.16%  '

# Row 10
Set-TextCell 'D10' '21.76'

# Row 11
Set-TextCell 'D11' '0.07808'
$ws.Range('E11').Value = '  +1.72%  '

# Row 12
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D12' '4.490'
$ws.Range('E12').Value = '  +1.49%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 'D13' '1.660.70'
$ws.Range('E13').Value = '  +0.87%  '

# Row 14
Set-TextCell 'D14' '0.5577'
$ws.Range('E14').Value = '  +0.61%  '

# Row 15
Set-TextCell 'D15' '0.0₅8308'
$ws.Range('E15').Value = '  +0.08%  '

# Row 16
Set-TextCell 'D16' '65.58'
$ws.Range('E16').Value = '  +1.26%  '

# Row 17
Set-TextCell 'D17' '26.494.86'

# Row 18
Set-TextCell 'D18' '1.001'
$ws.Range('E18').Value = '  -0.01%  '

# Row 19
Set-TextCell 'D19' '4.768'
$ws.Range('E19').Value = '  +1.15%  '

# Row 20
Set-TextCell 'D20' '193.22'
$ws.Range('E20').Value = '  +2.62%  '

# Row 21
$ws.Range('E21').Value = '  +1.51%  '

# Row 22
Set-TextCell 'D22' '6.311'
$ws.Range('E22').Value = '  +0.87%  '

# Row 23
$ws.Range('E23').Value = '  +0.00%  '

# Row 24
$ws.Range('E24').Value = '  +3.84%  '

# Row 25
Set-TextCell 'D25' '139.07'
$ws.Range('E25').Value = '  -3.87%  '

# Row 26
Set-TextCell 'D26' '7.394'
$ws.Range('E26').Value = '  -0.17%  '

# Row 27
Set-TextCell 'D27' '16.29'
$ws.Range('E27').Value = '  +2.57%  '

# Row 28
$ws.Range('E28').Value = '  +2.71%  '

# Row 29
Set-TextCell 'D29' '0.06215'
$ws.Range('E29').Value = '  +4.38%  '

# Row 30
Set-TextCell 'D30' '1.290'
$ws.Range('E30').Value = '  +1.98%  '

# Row 31
Set-TextCell 'D31' '3.608'
$ws.Range('E31').Value = '  +6.01%  '

# Row 32
Set-TextCell 'D32' '3.430'
$ws.Range('E32').Value = '  +0.91%  '

# Row 33
Set-TextCell 'D33' '1.685'
$ws.Range('E33').Value = '  +2.01%  '

# Row 34
$ws.Range('E34').Value = '  +1.13%  '

# Row 35
Set-TextCell 'D35' '0.6105'
$ws.Range('E35').Value = '  +8.37%  '

# Row 36
Set-TextCell 'D36' '2.413'
$ws.Range('E36').Value = '  +0.89%  '

# Row 37
Set-TextCell 'D37' '2.779'
$ws.Range('E37').Value = '  +0.97%  '

# Row 38
Set-TextCell 'D38' '0.01615'
$ws.Range('E38').Value = '  +0.44%  '

# Row 39
Set-TextCell 'D39' '6.041'
$ws.Range('E39').Value = '  +3.05%  '

# Row 40
Set-TextCell 'D40' '1.087.77'
$ws.Range('E40').Value = '  +5.83%  '

# Row 41
Set-TextCell 'D41' '0.8554'
$ws.Range('E41').Value = '  +0.35%  '

# Row 42
Set-TextCell 'D42' '1.000'
$ws.Range('E42').Value = '  -0.04%  '

# Row 43
Set-TextCell 'D43' '100.53'
$ws.Range('E43').Value = '  +1.75%  '

# Row 44
Set-TextCell 'D44' '1.816.69'
$ws.Range('E44').Value = '  +1.15%  '

# Row 45
$ws.Range('E45').Value = '  +3.78%  '

# Row 46
Set-TextCell 'D46' '58.34'
$ws.Range('E46').Value = '  +4.69%  '

# Row 47
$ws.Range('E47').Value = '  -0.58%  '

# Row 48
Set-TextCell 'D48' '1.521'
$ws.Range('E48').Value = '  +10.19%  '

# Row 49
Set-TextCell 'D49' '8.097'
$ws.Range('E49').Value = '  +0.48%  '

# Row 50
Set-TextCell 'D50' '0.05195'
$ws.Range('E50').Value = '  +1.00%  '

# Row 51
Set-TextCell 'D51' '6.002'
$ws.Range('E51').Value = '  +1.46%  '
